$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-16 15:35:18"
$wsZhCn.Range("G2").Value = "2016-02-16 15:36:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-16 15:35:37"
$wsDeDe.Range("G2").Value = "2016-02-16 15:36:44"
